$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Disponibilita" (availability) column C for rows 4-10 from 0 to 1
$ws.Range("C4:C10").Value = 1

# Update the active cell selection to C11
$ws.Range("C11").Select()
